# Weekly cryptos list refresh: updated Price/Volume(1h) figures for every coin row,
# and a handful of adjacent rows swapped which coin occupies that rank (Coin/Link/Price
# all change together there) -- e.g. Toncoin now ranks above InjectiveProtocol, etc.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "43.526.34"
$ws.Range("E2").Value = "  -6.38%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "2.593.03"
$ws.Range("E3").Value = "  -0.23%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.08%  "

# Row 5 (BNB)
$ws.Range("D5").Value = "'301.10"
$ws.Range("E5").Value = "  -2.35%  "

# Row 6 (Solana)
$ws.Range("D6").Value = "'96.33"
$ws.Range("E6").Value = "  -4.09%  "

# Row 7 (XRP)
$ws.Range("D7").Value = "'0.578"
$ws.Range("E7").Value = "  -4.40%  "

# Row 8 (USDC)
$ws.Range("E8").Value = "  +0.21%  "

# Row 9 (Cardano)
$ws.Range("D9").Value = "'0.559"
$ws.Range("E9").Value = "  -3.66%  "

# Row 10 (Avalanche)
$ws.Range("D10").Value = "'36.98"
$ws.Range("E10").Value = "  -6.08%  "

# Row 11 (Dogecoin)
$ws.Range("D11").Value = "'0.0815"
$ws.Range("E11").Value = "  -3.70%  "

# Row 12 (Polkadot)
$ws.Range("D12").Value = "'7.83"
$ws.Range("E12").Value = "  -4.20%  "

# Row 13 (WrappedliquidstakedEther2.0)
$ws.Range("D13").Value = "2.992.60"
$ws.Range("E13").Value = "  -0.05%  "

# Row 14 (TRON)
$ws.Range("E14").Value = "  +1.05%  "

# Row 15 (WrappedEther)
$ws.Range("D15").Value = "2.598.21"
$ws.Range("E15").Value = "  -0.16%  "

# Row 16 (Polygon)
$ws.Range("D16").Value = "'0.894"
$ws.Range("E16").Value = "  -3.33%  "

# Row 17 (Chainlink)
$ws.Range("D17").Value = "'14.39"
$ws.Range("E17").Value = "  -4.37%  "

# Row 18 (WrappedBTC)
$ws.Range("D18").Value = "43.635.22"
$ws.Range("E18").Value = "  -6.31%  "

# Row 19 (Uniswap)
$ws.Range("D19").Value = "'6.67"
$ws.Range("E19").Value = "  -1.63%  "

# Row 20 (ShibaInu)
$ws.Range("D20").Value = "0.0₃0979"
$ws.Range("E20").Value = "  -3.86%  "

# Row 21 (InternetComputer(DFINITY))
$ws.Range("D21").Value = "'12.37"
$ws.Range("E21").Value = "  -5.20%  "

# Row 22 (Litecoin)
$ws.Range("D22").Value = "'73.31"
$ws.Range("E22").Value = "  +2.05%  "

# Row 23 (BitcoinCash)
$ws.Range("D23").Value = "'266.72"
$ws.Range("E23").Value = "  -3.97%  "

# Row 24 (ImmutableX)
$ws.Range("E24").Value = "  +2.00%  "

# Row 25 (PancakeSwap)
$ws.Range("E25").Value = "  -3.88%  "

# Row 26 (EthereumClassic)
$ws.Range("D26").Value = "'29.49"
$ws.Range("E26").Value = "  +0.57%  "

# Row 27 (Dai)
$ws.Range("E27").Value = "  -0.11%  "

# Row 28 (Cosmos)
$ws.Range("D28").Value = "'10.29"
$ws.Range("E28").Value = "  -3.51%  "

# Row 29 (Toncoin)
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.22"
$ws.Range("E29").Value = "  -1.88%  "

# Row 30 (InjectiveProtocol)
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "'37.57"
$ws.Range("E30").Value = "  -3.97%  "

# Row 31 (Filecoin)
$ws.Range("D31").Value = "'6.05"
$ws.Range("E31").Value = "  -4.52%  "

# Row 32 (LidoDAOToken)
$ws.Range("E32").Value = "  +0.22%  "

# Row 33 (ARBITRUM)
$ws.Range("E33").Value = "  +2.33%  "

# Row 34 (Monero)
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").Value = "'152.07"
$ws.Range("E34").Value = "  +0.30%  "

# Row 35 (WEMIXToken)
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'2.79"
$ws.Range("E35").Value = "  -1.75%  "

# Row 36 (Hedera)
$ws.Range("D36").Value = "'0.0816"
$ws.Range("E36").Value = "  -2.87%  "

# Row 37 (Kaspa)
$ws.Range("E37").Value = "  -4.90%  "

# Row 38 (EnergySwap)
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").Value = "'24.51"
$ws.Range("E38").Value = "  +5.36%  "

# Row 39 (Stellar)
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "'0.121"
$ws.Range("E39").Value = "  -1.45%  "

# Row 40 (Celestia)
$ws.Range("D40").Value = "'17.04"
$ws.Range("E40").Value = "  +5.41%  "

# Row 41 (NEARProtocol)
$ws.Range("D41").Value = "'3.54"
$ws.Range("E41").Value = "  -3.14%  "

# Row 42 (VeChain)
$ws.Range("E42").Value = "  -5.33%  "

# Row 43 (RenderToken)
$ws.Range("D43").Value = "'3.85"
$ws.Range("E43").Value = "  -5.89%  "

# Row 44 (Maker)
$ws.Range("D44").Value = "2.067.96"
$ws.Range("E44").Value = "  -3.51%  "

# Row 45 (FirstDigitalUSD)
$ws.Range("D45").Value = "'0.998"
$ws.Range("E45").Value = "  +0.03%  "

# Row 46 (BitcoinSV)
$ws.Range("D46").Value = "'88.50"
$ws.Range("E46").Value = "  -4.99%  "

# Row 47 (FraxShare)
$ws.Range("D47").Value = "'9.07"
$ws.Range("E47").Value = "  -4.39%  "

# Row 48 (RocketPoolETH)
$ws.Range("D48").Value = "2.848.09"
$ws.Range("E48").Value = "  +0.08%  "

# Row 49 (ApeXProtocol)
$ws.Range("E49").Value = "  +1.04%  "

# Row 50 (Aave)
$ws.Range("D50").Value = "'106.01"
$ws.Range("E50").Value = "  -3.09%  "
